# Apply the cryptos-list price/volume refresh described in the commit diff.
# Source data is inline text (prices keep trailing zeros, thousands separated by
# dots, percentages padded with spaces) so every write must land as literal text,
# never get reinterpreted as a number. Cells whose new value parses as a plain
# decimal (e.g. "0.9998", "90.30") are written with a leading apostrophe, which
# is the standard COM/UI way to force text-entry without touching NumberFormat.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.762.54'
$ws.Range('E2').Value = '  -2.27%  '
$ws.Range('D3').Value = '1.796.03'
$ws.Range('D4').Value = '''0.9998'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''308.13'
$ws.Range('E5').Value = '  -1.81%  '
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').Value = '''0.4571'
$ws.Range('E7').Value = '  +1.79%  '
$ws.Range('E8').Value = '  -1.72%  '
$ws.Range('D9').Value = '''0.07242'
$ws.Range('E9').Value = '  -3.42%  '
$ws.Range('D10').Value = '''0.8549'
$ws.Range('E10').Value = '  -4.65%  '
$ws.Range('D11').Value = '''20.37'
$ws.Range('E11').Value = '  -3.51%  '
$ws.Range('D12').Value = '1.794.80'
$ws.Range('E12').Value = '  -1.76%  '
$ws.Range('D13').Value = '''5.304'
$ws.Range('E13').Value = '  -2.16%  '
$ws.Range('D14').Value = '''6.499'
$ws.Range('E14').Value = '  -4.15%  '
$ws.Range('D15').Value = '''0.07039'
$ws.Range('E15').Value = '  -1.21%  '
$ws.Range('D16').Value = '''90.30'
$ws.Range('E16').Value = '  -4.67%  '
$ws.Range('D17').Value = '''1.001'
$ws.Range('E17').Value = '  +0.07%  '
$ws.Range('D18').Value = '''0.000008628'
$ws.Range('E18').Value = '  -2.34%  '
$ws.Range('D19').Value = '''0.9992'
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('E20').Value = '  -4.16%  '
$ws.Range('D21').Value = '26.774.52'
$ws.Range('E21').Value = '  -2.28%  '
$ws.Range('D22').Value = '''5.285'
$ws.Range('E22').Value = '  -0.30%  '
$ws.Range('D23').Value = '''10.62'
$ws.Range('E23').Value = '  -3.44%  '
$ws.Range('D24').Value = '2.014.92'
$ws.Range('E24').Value = '  -1.77%  '
$ws.Range('D25').Value = '''1.907'
$ws.Range('E25').Value = '  -4.76%  '
$ws.Range('D26').Value = '''149.37'
$ws.Range('E26').Value = '  -1.72%  '
$ws.Range('D27').Value = '''2.160'
$ws.Range('E27').Value = '  -13.11%  '
$ws.Range('D28').Value = '''18.21'
$ws.Range('E28').Value = '  -2.33%  '
$ws.Range('D29').Value = '''5.200'
$ws.Range('E29').Value = '  -3.64%  '
$ws.Range('D30').Value = '''114.09'
$ws.Range('E30').Value = '  -3.55%  '
$ws.Range('D31').Value = '''0.08836'
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('D32').Value = '''0.7573'
$ws.Range('E32').Value = '  -3.25%  '
$ws.Range('D33').Value = '''1.159'
$ws.Range('E33').Value = '  -3.54%  '
$ws.Range('D34').Value = '''4.439'
$ws.Range('E34').Value = '  -3.38%  '
$ws.Range('D35').Value = '''2.880'
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('D36').Value = '''0.9990'
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('D37').Value = '''1.112'
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('D38').Value = '''0.01940'
$ws.Range('E38').Value = '  -2.67%  '
$ws.Range('D39').Value = '''0.05214'
$ws.Range('E39').Value = '  -2.45%  '
$ws.Range('D40').Value = '''2.372'
$ws.Range('E40').Value = '  +3.93%  '
$ws.Range('D41').Value = '''2.902'
$ws.Range('E41').Value = '  +1.56%  '
$ws.Range('D42').Value = '''7.131'
$ws.Range('E42').Value = '  -3.99%  '
$ws.Range('D43').Value = '''0.5241'
$ws.Range('E43').Value = '  -2.24%  '
$ws.Range('D44').Value = '''0.1646'
$ws.Range('E44').Value = '  -5.16%  '
$ws.Range('D45').Value = '''0.5074'
$ws.Range('E45').Value = '  -1.66%  '
$ws.Range('D46').Value = '''8.481'
$ws.Range('E46').Value = '  -4.17%  '
$ws.Range('D47').Value = '''10.21'
$ws.Range('E47').Value = '  -5.24%  '
$ws.Range('D48').Value = '''104.24'
$ws.Range('E48').Value = '  -2.39%  '
$ws.Range('D49').Value = '''0.9981'
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('D50').Value = '''1.647'
$ws.Range('E50').Value = '  -3.71%  '
$ws.Range('D51').Value = '''0.06295'
$ws.Range('E51').Value = '  -1.49%  '
